# Add new vehicle data rows to the worksheet (rows 13-16), which also
# appends the corresponding new unique strings to the shared-strings table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Furgone / fiat / doblo / 2020
$ws.Range("B13").Value = "Furgone"
$ws.Range("D13").Value = "fiat"
$ws.Range("E13").Value = "doblo"
$ws.Range("F13").Value = 2020

# Row 14: Furgone / fiat / doblò / 2020
$ws.Range("B14").Value = "Furgone"
$ws.Range("D14").Value = "fiat"
$ws.Range("E14").Value = "doblò"
$ws.Range("F14").Value = 2020

# Row 15: Furgone / fiat / doblo / 2017
$ws.Range("B15").Value = "Furgone"
$ws.Range("D15").Value = "fiat"
$ws.Range("E15").Value = "doblo"
$ws.Range("F15").Value = 2017

# Row 16: adsf / Veicolo / bmw / 1-series / 2017
$ws.Range("A16").Value = "adsf"
$ws.Range("B16").Value = "Veicolo"
$ws.Range("D16").Value = "bmw"
$ws.Range("E16").Value = "1-series"
$ws.Range("F16").Value = 2017

# Update the active selection to match the cell left selected after entry
[void]$ws.Range("H17").Select()
